# Apply edits to parallel.xlsx data: extend grid from O to Q (add columns P,Q),
# swap values in columns I/K/M/O for rows 2-25, and add new header values in P1/Q1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: add P1 = 14, Q1 = 15, matching the style of O1 (index header style) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Copy formatting (bold, centered, border) from O1 to the new header cells
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats = -4122
$excel.CutCopyMode = 0

# --- Data rows 2-25: swap I/K/M/O values and add P/Q = 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new = 2
}

$wb.Save()
